$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date column (A2:A5) to the new reporting date (2023-09-01)
$ws.Range("A2:A5").Value = [DateTime]"2023-09-01"

# Update the quantities (column C) for the new period
$ws.Range("C2").Value = 1457
$ws.Range("C3").Value = 103
$ws.Range("C4").Value = 321
$ws.Range("C5").Value = 77

# Move the active selection to match the saved view state
$ws.Range("F10").Select()
